$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.762.68"
$ws.Range("E2").Value = "  +0.77%  "

$ws.Range("D3").Value = "1.853.66"
$ws.Range("E3").Value = "  +0.60%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.86"
$ws.Range("E5").Value = "  -0.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6553"
$ws.Range("E6").Value = "  +4.22%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.06"
$ws.Range("E8").Value = "  +3.95%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07511"
$ws.Range("E9").Value = "  +0.70%  "

$ws.Range("E10").Value = "  +0.62%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "24.58"
$ws.Range("E11").Value = "  +3.51%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07641"
$ws.Range("E12").Value = "  -0.36%  "

$ws.Range("D13").Value = "1.858.97"
$ws.Range("E13").Value = "  +1.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.052"
$ws.Range("E14").Value = "  +0.38%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6865"
$ws.Range("E15").Value = "  +0.95%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.63"
$ws.Range("E16").Value = "  -0.69%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000009702"
$ws.Range("E17").Value = "  +3.50%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.109"
$ws.Range("E18").Value = "  +2.68%  "

$ws.Range("D19").Value = "29.783.28"
$ws.Range("E19").Value = "  +0.92%  "

$ws.Range("D20").Value = "2.108.00"
$ws.Range("E20").Value = "  +1.18%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "237.02"
$ws.Range("E21").Value = "  -0.27%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.64"
$ws.Range("E22").Value = "  +0.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9998"
$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.722"
$ws.Range("E24").Value = "  +4.64%  "

$ws.Range("E25").Value = "  +0.11%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.27"
$ws.Range("E26").Value = "  -0.69%  "

$ws.Range("E27").Value = "  +1.03%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.541"
$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "17.86"
$ws.Range("E29").Value = "  +0.24%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06069"
$ws.Range("E30").Value = "  -0.51%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.491"
$ws.Range("E31").Value = "  -0.32%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.275"
$ws.Range("E32").Value = "  +2.55%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.142"
$ws.Range("E33").Value = "  +0.38%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.076"
$ws.Range("E34").Value = "  -0.71%  "

$ws.Range("E35").Value = "  -0.14%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.181"
$ws.Range("E36").Value = "  +2.95%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.7249"
$ws.Range("E37").Value = "  -0.64%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.607"
$ws.Range("E38").Value = "  +0.08%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.804"
$ws.Range("E39").Value = "  -2.53%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01793"
$ws.Range("E40").Value = "  +1.63%  "

$ws.Range("D41").Value = "1.200.54"
$ws.Range("E41").Value = "  -1.97%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.285"
$ws.Range("E42").Value = "  -0.19%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9115"
$ws.Range("E43").Value = "  -0.37%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9998"
$ws.Range("E44").Value = "  -0.09%  "

$ws.Range("D45").Value = "2.011.90"
$ws.Range("E45").Value = "  +0.70%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.20"
$ws.Range("E46").Value = "  -0.83%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "66.80"
$ws.Range("E47").Value = "  +1.63%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.331"
$ws.Range("E48").Value = "  +9.70%  "

$ws.Range("E49").Value = "  +2.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4061"
$ws.Range("E50").Value = "  -0.22%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.111"
$ws.Range("E51").Value = "  -1.46%  "
